$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap between row 46 and row 47 (everything except A=row# and E=date, which are
# identical in nature per-row / already correct): B (id), F (HomeTeam), G (AwayTeam),
# H (FTHG), I (FTAG), J (FTR), K..AC (odds)
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($col in $cols) {
    $addr46 = "$col" + "46"
    $addr47 = "$col" + "47"
    $val46 = $ws.Range($addr46).Value2
    $val47 = $ws.Range($addr47).Value2
    $ws.Range($addr46).Value = $val47
    $ws.Range($addr47).Value = $val46
}
